# Update TPM-derived NATMI edge metrics for Tnf-Tnfrsf21 (YoungD7).
#
# The Ligand average/total expression columns (G/H) are keyed only by the
# Sending cluster (col A); the Receptor average/total expression columns
# (M/N) are keyed only by the Target cluster (col D). Refreshing the TPM
# pipeline changed the per-cluster ligand values for ECs and
# Resolving-Mac, and the per-cluster receptor values for ECs, MuSCs and
# Resolving-Mac (FAPs stayed the same on both sides). Every other changed
# column is derived from those:
#   I = G / sum(G over clusters)      J = H / sum(H over clusters)
#   O = M / sum(M over clusters)      P = N / sum(N over clusters)
#   Q = G * M                         R = H * N
#   S = Q / sum(Q over all rows)      T = R / sum(R over all rows)
# so this script simply writes the refreshed values straight onto the
# affected cells, matching the new TPM output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 2.624526
$ws.Range("H2").Value2 = 7.873578
$ws.Range("I2").Value2 = 0.06442870872595916
$ws.Range("J2").Value2 = 0.06442870872595916
$ws.Range("M2").Value2 = 8.461686666666667
$ws.Range("N2").Value2 = 25.38506
$ws.Range("O2").Value2 = 0.1873739652872041
$ws.Range("P2").Value2 = 0.1873739652872041
$ws.Range("Q2").Value2 = 22.20791666052
$ws.Range("R2").Value2 = 199.87124994468
$ws.Range("S2").Value2 = 0.01207226263231726
$ws.Range("T2").Value2 = 0.01207226263231726
$ws.Range("G3").Value2 = 2.624526
$ws.Range("H3").Value2 = 7.873578
$ws.Range("I3").Value2 = 0.06442870872595916
$ws.Range("J3").Value2 = 0.06442870872595916
$ws.Range("O3").Value2 = 0.02691675086296081
$ws.Range("P3").Value2 = 0.02691675086296081
$ws.Range("Q3").Value2 = 3.190224207618
$ws.Range("R3").Value2 = 28.712017868562
$ws.Range("S3").Value2 = 0.001734211501198912
$ws.Range("T3").Value2 = 0.001734211501198912
$ws.Range("G4").Value2 = 2.624526
$ws.Range("H4").Value2 = 7.873578
$ws.Range("I4").Value2 = 0.06442870872595916
$ws.Range("J4").Value2 = 0.06442870872595916
$ws.Range("M4").Value2 = 16.89366666666666
$ws.Range("N4").Value2 = 50.681
$ws.Range("O4").Value2 = 0.3740901118500721
$ws.Range("P4").Value2 = 0.3740901118500721
$ws.Range("Q4").Value2 = 44.33786740199999
$ws.Range("R4").Value2 = 399.040806618
$ws.Range("S4").Value2 = 0.02410214285364978
$ws.Range("T4").Value2 = 0.02410214285364978
$ws.Range("G5").Value2 = 2.624526
$ws.Range("H5").Value2 = 7.873578
$ws.Range("I5").Value2 = 0.06442870872595916
$ws.Range("J5").Value2 = 0.06442870872595916
$ws.Range("M5").Value2 = 18.58845466666667
$ws.Range("N5").Value2 = 55.76536400000001
$ws.Range("O5").Value2 = 0.4116191719997629
$ws.Range("P5").Value2 = 0.411619171999763
$ws.Range("Q5").Value2 = 48.785882572488
$ws.Range("R5").Value2 = 439.072943152392
$ws.Range("S5").Value2 = 0.02652009173879321
$ws.Range("T5").Value2 = 0.02652009173879321
$ws.Range("I6").Value2 = 0.01210207677934112
$ws.Range("J6").Value2 = 0.01210207677934112
$ws.Range("M6").Value2 = 8.461686666666667
$ws.Range("N6").Value2 = 25.38506
$ws.Range("O6").Value2 = 0.1873739652872041
$ws.Range("P6").Value2 = 0.1873739652872041
$ws.Range("Q6").Value2 = 4.171462036868889
$ws.Range("R6").Value2 = 37.54315833182
$ws.Range("S6").Value2 = 0.002267614114355343
$ws.Range("T6").Value2 = 0.002267614114355343
$ws.Range("I7").Value2 = 0.01210207677934112
$ws.Range("J7").Value2 = 0.01210207677934112
$ws.Range("O7").Value2 = 0.02691675086296081
$ws.Range("P7").Value2 = 0.02691675086296081
$ws.Range("S7").Value2 = 0.0003257485855939481
$ws.Range("T7").Value2 = 0.0003257485855939482
$ws.Range("I8").Value2 = 0.01210207677934112
$ws.Range("J8").Value2 = 0.01210207677934112
$ws.Range("M8").Value2 = 16.89366666666666
$ws.Range("N8").Value2 = 50.681
$ws.Range("O8").Value2 = 0.3740901118500721
$ws.Range("P8").Value2 = 0.3740901118500721
$ws.Range("Q8").Value2 = 8.328279211888889
$ws.Range("R8").Value2 = 74.954512907
$ws.Range("S8").Value2 = 0.004527267256001881
$ws.Range("T8").Value2 = 0.004527267256001881
$ws.Range("I9").Value2 = 0.01210207677934112
$ws.Range("J9").Value2 = 0.01210207677934112
$ws.Range("M9").Value2 = 18.58845466666667
$ws.Range("N9").Value2 = 55.76536400000001
$ws.Range("O9").Value2 = 0.4116191719997629
$ws.Range("P9").Value2 = 0.411619171999763
$ws.Range("Q9").Value2 = 9.163779754634223
$ws.Range("R9").Value2 = 82.47401779170801
$ws.Range("S9").Value2 = 0.00498144682338995
$ws.Range("T9").Value2 = 0.004981446823389951
$ws.Range("I10").Value2 = 0.003429134645952472
$ws.Range("J10").Value2 = 0.003429134645952472
$ws.Range("M10").Value2 = 8.461686666666667
$ws.Range("N10").Value2 = 25.38506
$ws.Range("O10").Value2 = 0.1873739652872041
$ws.Range("P10").Value2 = 0.1873739652872041
$ws.Range("Q10").Value2 = 1.181987625406667
$ws.Range("R10").Value2 = 10.63788862866
$ws.Range("S10").Value2 = 0.0006425305561158474
$ws.Range("T10").Value2 = 0.0006425305561158474
$ws.Range("I11").Value2 = 0.003429134645952472
$ws.Range("J11").Value2 = 0.003429134645952472
$ws.Range("O11").Value2 = 0.02691675086296081
$ws.Range("P11").Value2 = 0.02691675086296081
$ws.Range("S11").Value2 = 0.00009230116294064998
$ws.Range("T11").Value2 = 0.00009230116294065
$ws.Range("I12").Value2 = 0.003429134645952472
$ws.Range("J12").Value2 = 0.003429134645952472
$ws.Range("M12").Value2 = 16.89366666666666
$ws.Range("N12").Value2 = 50.681
$ws.Range("O12").Value2 = 0.3740901118500721
$ws.Range("P12").Value2 = 0.3740901118500721
$ws.Range("Q12").Value2 = 2.359825615666666
$ws.Range("R12").Value2 = 21.238430541
$ws.Range("S12").Value2 = 0.001282805363253317
$ws.Range("T12").Value2 = 0.001282805363253318
$ws.Range("I13").Value2 = 0.003429134645952472
$ws.Range("J13").Value2 = 0.003429134645952472
$ws.Range("M13").Value2 = 18.58845466666667
$ws.Range("N13").Value2 = 55.76536400000001
$ws.Range("O13").Value2 = 0.4116191719997629
$ws.Range("P13").Value2 = 0.4116191719997629
$ws.Range("Q13").Value2 = 2.596565467022667
$ws.Range("R13").Value2 = 23.369089203204
$ws.Range("S13").Value2 = 0.001411497563642657
$ws.Range("T13").Value2 = 0.001411497563642657
$ws.Range("G14").Value2 = 37.47815466666666
$ws.Range("H14").Value2 = 112.434464
$ws.Range("I14").Value2 = 0.9200400798487472
$ws.Range("J14").Value2 = 0.9200400798487472
$ws.Range("M14").Value2 = 8.461686666666667
$ws.Range("N14").Value2 = 25.38506
$ws.Range("O14").Value2 = 0.1873739652872041
$ws.Range("P14").Value2 = 0.1873739652872041
$ws.Range("Q14").Value2 = 317.1284016342044
$ws.Range("R14").Value2 = 2854.15561470784
$ws.Range("S14").Value2 = 0.1723915579844157
$ws.Range("T14").Value2 = 0.1723915579844157
$ws.Range("G15").Value2 = 37.47815466666666
$ws.Range("H15").Value2 = 112.434464
$ws.Range("I15").Value2 = 0.9200400798487472
$ws.Range("J15").Value2 = 0.9200400798487472
$ws.Range("O15").Value2 = 0.02691675086296081
$ws.Range("P15").Value2 = 0.02691675086296081
$ws.Range("Q15").Value2 = 45.556308557984
$ws.Range("R15").Value2 = 410.006777021856
$ws.Range("S15").Value2 = 0.0247644896132273
$ws.Range("T15").Value2 = 0.0247644896132273
$ws.Range("G16").Value2 = 37.47815466666666
$ws.Range("H16").Value2 = 112.434464
$ws.Range("I16").Value2 = 0.9200400798487472
$ws.Range("J16").Value2 = 0.9200400798487472
$ws.Range("M16").Value2 = 16.89366666666666
$ws.Range("N16").Value2 = 50.681
$ws.Range("O16").Value2 = 0.3740901118500721
$ws.Range("P16").Value2 = 0.3740901118500721
$ws.Range("Q16").Value2 = 633.1434522204443
$ws.Range("R16").Value2 = 5698.291069983999
$ws.Range("S16").Value2 = 0.3441778963771671
$ws.Range("T16").Value2 = 0.3441778963771672
$ws.Range("G17").Value2 = 37.47815466666666
$ws.Range("H17").Value2 = 112.434464
$ws.Range("I17").Value2 = 0.9200400798487472
$ws.Range("J17").Value2 = 0.9200400798487472
$ws.Range("M17").Value2 = 18.58845466666667
$ws.Range("N17").Value2 = 55.76536400000001
$ws.Range("O17").Value2 = 0.4116191719997629
$ws.Range("P17").Value2 = 0.411619171999763
$ws.Range("Q17").Value2 = 696.6609790116551
$ws.Range("R17").Value2 = 6269.948811104896
$ws.Range("S17").Value2 = 0.3787061358739371
$ws.Range("T17").Value2 = 0.3787061358739372
